# Add support for Trillion units
# Adds a new "IncomeEquities" worksheet with Disposable Personal Income
# and Corporate Equities & Mutual Fund Shares data (quarterly, 2023-2025).

$wb = $excel.ActiveWorkbook

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "IncomeEquities"

# Move the new sheet to the end of the tab strip (after the last
# pre-existing sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# Re-acquire the worksheet reference by name, since the in-place Move
# invalidates the old handle's positional binding.
$ws = $wb.Worksheets.Item("IncomeEquities")

$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Disposable Personal Income"
$ws.Range("C1").Value = "Corporate Equities & Mutual Fund Shares"

$dates = @("1/1/2023","4/1/2023","7/1/2023","10/1/2023","1/1/2024","4/1/2024","7/1/2024","10/1/2024","1/1/2025","4/1/2025")
$dpi = @(20283400000000,20651000000000,20894600000000,21168000000000,21575400000000,21843200000000,22002600000000,22249500000000,22563700000000,22858500000000)
$ceq = @(36520700000000,38627200000000,37204700000000,40738200000000,44267800000000,45057300000000,48147900000000,48547800000000,46723700000000,51186500000000)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $dpi[$i]
    $ws.Cells.Item($row, 3).Value = $ceq[$i]
}

$ws.Activate()
